$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)
$shp = $s.Shapes.AddShape(1, 2254102/12700, 6188149/12700, 1765005/12700, 489098/12700)
$shp.Name = "Rectangle 1"
